# Trade #17 closed at 2026-02-17 20:52:56 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: roll up totals after the new trade opens + trade #45 closes.
#  - Strategy Status sheet: MarketMaking strategy row totals.
#  - All Trades sheet: trade #45 (row 46) closes out; trade #78 (row 79) opens.
#  - MarketMaking sheet: trade #45 (row 13) closes out; trade #78 (row 46) opens.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.4    # Current Capital
$summary.Range("B4").Value = 0.19      # Total P&L $
$summary.Range("B5").Value = 0.08      # Total P&L %
$summary.Range("B6").Value = 45        # Total Trades
$summary.Range("B7").Value = 20        # Winning Trades
$summary.Range("B9").Value = 44.44     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.4      # Capital
$status.Range("D5").Value = 12         # Trades
$status.Range("E5").Value = 0.08       # P&L $
$status.Range("F5").Value = 0.4        # P&L %
$status.Range("G5").Value = 50         # Win Rate %

# ---------------------------------------------------------------------------
# All Trades
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #45 (row 46) transitions from OPEN to CLOSED.
$allTrades.Range("G46").Value = 0.878646
$allTrades.Range("H46").Value = "CLOSED"
$allTrades.Range("I46").Value = 5.8609
$allTrades.Range("J46").Value = 0.05
$allTrades.Range("K46").Value = 100.4
$allTrades.Range("L46").Value = "early_exit"
$allTrades.Range("M46").Value = 0.14

# Trade #78 (row 79) is newly appended / opened.
$allTrades.Range("A79").Value = 78
$allTrades.Range("B79").NumberFormat = "@"
$allTrades.Range("B79").Value = "2026-02-17"
$allTrades.Range("C79").Value = "20:52:50"
$allTrades.Range("D79").Value = "MarketMaking"
$allTrades.Range("E79").Value = "DOWN"
$allTrades.Range("F79").Value = 0.83
$allTrades.Range("H79").Value = "OPEN"
$allTrades.Range("I79").Value = 0
$allTrades.Range("J79").Value = 0
$allTrades.Range("K79").Value = 100.3523945789973
$allTrades.Range("M79").Value = 0
$allTrades.Range("N79").Value = 0
$allTrades.Range("O79").Value = 0
$allTrades.Range("P79").Value = 0.6
$allTrades.Range("Q79").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Trade #45 (row 13) transitions from OPEN to CLOSED.
$mm.Range("G13").Value = 0.878646
$mm.Range("H13").Value = "CLOSED"
$mm.Range("I13").Value = 5.8609
$mm.Range("J13").Value = 0.05
$mm.Range("K13").Value = 100.4
$mm.Range("P13").Value = "early_exit"
$mm.Range("Q13").Value = 0.14

# Trade #78 (row 46) is newly appended / opened.
$mm.Range("A46").Value = 78
$mm.Range("B46").NumberFormat = "@"
$mm.Range("B46").Value = "2026-02-17"
$mm.Range("C46").Value = "20:52:50"
$mm.Range("D46").Value = "MarketMaking"
$mm.Range("E46").Value = "DOWN"
$mm.Range("F46").Value = 0.83
$mm.Range("H46").Value = "OPEN"
$mm.Range("I46").Value = 0
$mm.Range("J46").Value = 0
$mm.Range("K46").Value = 100.3523945789973
$mm.Range("L46").Value = 0
$mm.Range("M46").Value = 0
$mm.Range("N46").Value = 0.6
$mm.Range("O46").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q46").Value = 0
